# Coding_Tracker.xlsx edit:
#   Insert a new "Greedy" worksheet between "Hashing" and "Dynamic Programming",
#   populate it with 12 Greedy-algorithm problems (first marked "Done", with a
#   hyperlink to its GfG page, plus hyperlinks for the remaining 11 rows), and
#   leave the surrounding sheets' selection/active-tab state the way Excel
#   would after that edit.

$wb = $excel.ActiveWorkbook

# --- 1. Stash the "previous" selection on sheets whose stored cursor moves
#        as a side effect of this edit (captured BEFORE the new sheet steals
#        the active-tab flag). Doing the Range.Select() calls up front, then
#        re-selecting the final cell on the new sheet last, reproduces what
#        Excel itself persists into each sheetView.
$hashing = $wb.Worksheets.Item("Hashing")
$hashing.Range("A1:XFD1").Select()

$dp = $wb.Worksheets.Item("Dynamic Programming")
$dp.Range("A2").Select()

# --- 2. Insert the new "Greedy" sheet right after "Hashing" -----------------
$greedy = $wb.Worksheets.Add([Type]::Missing, $hashing)
$greedy.Name = "Greedy"

# --- 3. Header row ------------------------------------------------------
$greedy.Range("A1").Value = "Problem"
$greedy.Range("B1").Value = "Status (Done/Not Done)"

# --- 4. Problem rows ------------------------------------------------------
$problems = @(
    @{Row=2;  Title="1. Activity Selection";              Status="Done"; Url="https://www.geeksforgeeks.org/problems/activity-selection-1587115620/1"},
    @{Row=3;  Title="2. N meetings in one room";           Status=$null;  Url="https://www.geeksforgeeks.org/problems/n-meetings-in-one-room-1587115620/1"},
    @{Row=4;  Title="3. Coin Piles";                       Status=$null;  Url="https://www.geeksforgeeks.org/problems/choose-and-swap0531/1"},
    @{Row=5;  Title="4. Maximize Toys";                    Status=$null;  Url="https://www.geeksforgeeks.org/problems/maximize-toys0331/1"},
    @{Row=6;  Title="5. Page Faults in LRU";                Status=$null;  Url="https://www.geeksforgeeks.org/problems/page-faults-in-lru5603/1"},
    @{Row=7;  Title="6. Largest number possible";           Status=$null;  Url="https://www.geeksforgeeks.org/problems/largest-number-possible5028/1"},
    @{Row=8;  Title="7. Minimize the heights";              Status=$null;  Url="https://www.geeksforgeeks.org/problems/minimize-the-heights3351/1"},
    @{Row=9;  Title="8. Minimize the sum of product";       Status=$null;  Url="https://www.geeksforgeeks.org/problems/minimize-the-sum-of-product1525/1"},
    @{Row=10; Title="9. Huffman Decoding";                  Status=$null;  Url="https://www.geeksforgeeks.org/problems/huffman-decoding-1/1"},
    @{Row=11; Title="10. Minimum Spanning Tree";            Status=$null;  Url="https://www.geeksforgeeks.org/problems/minimum-spanning-tree/1"},
    @{Row=12; Title="11. Shop in Candy Store";              Status=$null;  Url="https://www.geeksforgeeks.org/problems/shop-in-candy-store1145/1"},
    @{Row=13; Title="12. Geek collects the balls";          Status=$null;  Url="https://www.geeksforgeeks.org/problems/geek-collects-the-balls5515/1"}
)

foreach ($p in $problems) {
    $cell = $greedy.Cells.Item($p.Row, 1)
    $cell.Value = $p.Title
    $greedy.Hyperlinks.Add($cell, $p.Url, [Type]::Missing, [Type]::Missing, $p.Url) | Out-Null
    # Re-assert the plain text value/style - Excel's own "Activity Selection"
    # row (and every other Greedy row) carries no hyperlink styling in the
    # source sheet, so strip the auto-applied Hyperlink cell style back off.
    $cell.Value = $p.Title
    $cell.Style = "Normal"

    if ($p.Status) {
        $greedy.Cells.Item($p.Row, 2).Value = $p.Status
    }
}

# --- 5. Column widths (best effort bestFit) --------------------------------
$greedy.Columns.Item(1).AutoFit() | Out-Null
$greedy.Columns.Item(2).AutoFit() | Out-Null

# --- 6. Selection/active-tab bookkeeping -----------------------------------
$greedy.Range("B2").Select()
